# Update "horarios" workbook: Linea 141 schedules refreshed to 05:22:23.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 05:22:23"
$ws1.Range("A3").Value = "Total filas: 22"

$data1 = @(
    @("05:22:23", "05:22", "23_HERNANDEZ", 0, "LP1912"),
    @("05:22:23", "05:34", "215B_EL PATO", 12, "LP1912"),
    @("05:22:23", "05:46", "15_ABASTO", 24, "LP1912"),
    @("05:22:23", "05:54", "10_OLMOS", 32, "LP1912"),
    @("05:22:23", "06:04", "16_SANTA ANA", 42, "LP1912"),
    @("05:22:23", "06:11", "215A_EL PATO", 49, "LP1912"),
    @("05:22:23", "06:14", "225_HARAS DEL SUR", 52, "LP1912"),
    @("05:22:23", "06:21", "26_HERNANDEZ", 59, "LP1912"),
    @("05:22:23", "06:27", "23_HERNANDEZ", 65, "LP1912"),
    @("05:22:23", "06:29", "86_EST CHICA-ESC AGRARIA", 67, "LP1912"),
    @("05:22:23", "06:31", "16_SANTA ANA", 69, "LP1912"),
    @("05:22:23", "06:44", "225_C ROCA-H SUR", 82, "LP1912"),
    @("05:22:23", "06:46", "215C_EL PATO", 84, "LP1912"),
    @("05:22:23", "06:59", "14_ABASTO", 97, "LP1912"),
    @("05:22:23", "07:05", "15_ABASTO", 103, "LP1912"),
    @("05:22:23", "07:07", "225_GOMEZ", 105, "LP1912"),
    @("05:22:23", "07:11", "215A_EL PATO", 109, "LP1912"),
    @("05:22:23", "07:15", "11_ETCHEVERRY", 113, "LP1912"),
    @("05:22:23", "07:21", "26_HERNANDEZ", 119, "LP1912")
)

$startRow = 9
for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = $startRow + $i
    $row = $data1[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 05:22:23"
$ws2.Range("A3").Value = "Total filas: 5"

$data2 = @(
    @("05:22:23", "05:34", "215B_EL PATO", 12, "LP1912"),
    @("05:22:23", "06:11", "215A_EL PATO", 49, "LP1912"),
    @("05:22:23", "06:46", "215C_EL PATO", 84, "LP1912"),
    @("05:22:23", "07:11", "215A_EL PATO", 109, "LP1912")
)

$startRow2 = 7
for ($i = 0; $i -lt $data2.Length; $i++) {
    $r = $startRow2 + $i
    $row = $data2[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 05:22:23"
$ws3.Range("A3").Value = "Total filas: 7"

# Row 7 updates in place.
$ws3.Cells.Item(7, 1).Value = "05:22:23"
$ws3.Cells.Item(7, 2).Value = "05:44"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 22
$ws3.Cells.Item(7, 5).Value = "L6173"

# Row 8 is untouched (04:49:42 / 06:08 / 215A_LA PLATA / 79 / L6173).

# New rows 9-12 (old row 9 shifts down to row 10).
$ws3.Cells.Item(9, 1).Value = "05:22:23"
$ws3.Cells.Item(9, 2).Value = "06:09"
$ws3.Cells.Item(9, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(9, 4).Value = 47
$ws3.Cells.Item(9, 5).Value = "L6173"

$ws3.Cells.Item(10, 1).Value = "04:49:42"
$ws3.Cells.Item(10, 2).Value = "06:32"
$ws3.Cells.Item(10, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(10, 4).Value = 103
$ws3.Cells.Item(10, 5).Value = "L6203"

$ws3.Cells.Item(11, 1).Value = "05:22:23"
$ws3.Cells.Item(11, 2).Value = "06:33"
$ws3.Cells.Item(11, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(11, 4).Value = 71
$ws3.Cells.Item(11, 5).Value = "L6203"

$ws3.Cells.Item(12, 1).Value = "05:22:23"
$ws3.Cells.Item(12, 2).Value = "07:00"
$ws3.Cells.Item(12, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws3.Cells.Item(12, 4).Value = 98
$ws3.Cells.Item(12, 5).Value = "L6173"
